$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 89. This shifts the existing rows 89-106 down to 90-107,
# carrying their values/styles along with them (standard Excel row-insert behaviour).
$ws.Rows.Item(89).Insert()

# Populate the newly inserted row 89 with a new weekly record. All descriptive columns
# (A,B,C,E,F,G,H,I,J,K,Q,R,T) match the rest of this product's block; only the date,
# quality, prices and $/Kg differ for this new entry.
$ws.Cells.Item(89, 1).Value = 5
$ws.Cells.Item(89, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(89, 3).Value = "Maule"
$ws.Cells.Item(89, 4).Value = 44951
$ws.Cells.Item(89, 4).NumberFormat = $ws.Cells.Item(90, 4).NumberFormat
$ws.Cells.Item(89, 5).Value = 7
$ws.Cells.Item(89, 6).Value = "Fruta"
$ws.Cells.Item(89, 7).Value = 100101
$ws.Cells.Item(89, 8).Value = "Berries"
$ws.Cells.Item(89, 9).Value = 100101001
$ws.Cells.Item(89, 10).Value = "Arándano (blue)"
$ws.Cells.Item(89, 11).Value = "Sin especificar"
$ws.Cells.Item(89, 12).Value = "Segunda"
$ws.Cells.Item(89, 13).Value = 40
$ws.Cells.Item(89, 14).Value = 2600
$ws.Cells.Item(89, 15).Value = 2600
$ws.Cells.Item(89, 16).Value = 2600
$ws.Cells.Item(89, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(89, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(89, 19).Value = 1300
$ws.Cells.Item(89, 20).Value = 2
